$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (unchanged text, but shared string indices change internally - values stay same)
$ws.Range("A1").Value = "nis"
$ws.Range("B1").Value = "nama"
$ws.Range("C1").Value = "gambar"
$ws.Range("D1").Value = "quote"
$ws.Range("E1").Value = "jurusan"
$ws.Range("F1").Value = "d_kelas"

# Row 2
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "not error"
$ws.Range("C2").Value = "Si_A.jpg"
$ws.Range("D2").Value = "ah masa"
$ws.Range("E2").Value = "tei"
$ws.Range("F2").Value = 1

# Row 3
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "not error"
$ws.Range("C3").Value = "Si_B.jpg"
$ws.Range("D3").Value = "ah masa"
$ws.Range("E3").Value = "tei"
$ws.Range("F3").Value = 4

# Row 4 (new)
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "not error"
$ws.Range("C4").Value = "Si_C.jpg"
$ws.Range("D4").Value = "ah masa"
$ws.Range("E4").Value = "tbsm"
$ws.Range("F4").Value = 1

# Row 5 (new)
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "not error"
$ws.Range("C5").Value = "si_D.jpg"
$ws.Range("D5").Value = "ah masa"
$ws.Range("E5").Value = "TKJ"
$ws.Range("F5").Value = 4

# Update selection to F5 to mirror cursor position after edit
$ws.Range("F5").Select()
